$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "H2"  = 65
    "I2"  = 180
    "J2"  = 768
    "K2"  = 3
    "L2"  = 220
    "M2"  = 14
    "N2"  = 152
    "O2"  = 0
    "P2"  = 1
    "Q2"  = 0
    "R2"  = 18
    "S2"  = 71
    "T2"  = 152
    "U2"  = 11
    "V2"  = 1201
    "W2"  = 0
    "X2"  = 1198
    "Y2"  = 4
    "Z2"  = 15
    "AA2" = 6
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
